# Generate Report for Handback
# Update "Latest Handback DateTime" (column K) for row 2 (the
# 9a42bf13-1086-426d-9b7f-29691f2b0b05.md file) on both the zh-cn and
# de-de worksheets to reflect the freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-09-09 08:39:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-09 08:39:47"
